$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing data rows (2-48) down to (3-49).
$ws.Range("A2:R2").EntireRow.Insert()

# The inserted row picks up the header row's formatting by default; reset it
# back to the plain "Normal" style used by the other data rows, then restore
# the date number format on column D (matches the rest of the Fecha column).
$ws.Range("A2:R2").Style = "Normal"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the newly inserted row 2 with this week's data.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").Value = 44530
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 100112012
$ws.Range("G2").Value = "Espinaca"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 1900
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = 1950
$ws.Range("N2").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 650
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = "Hortaliza"
